# Adding the changes we made on may 9th
#
# This script:
#  1. Inserts 9 new data rows right after the header row (rows 2-10), pushing the
#     existing 20 data rows down to rows 11-30.
#  2. Fills the 9 newly inserted rows with new sensor readings.
#  3. Appends one brand new data row (row 31) with new sensor readings.
#  4. Recomputes column A (timestamp) for every data row as a contiguous
#     sequence 0, 100, 200, ... (100 per row), and makes sure column B
#     stays "falling" throughout.
#  5. Updates the used-range dimension to A1:H31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert 9 blank rows after the header (before the old row 2) ---
$ws.Rows("2:10").Insert()

# --- Step 2: new data for the 9 newly inserted rows (rows 2-10) ---
$newTopData = @(
    @(-2.896898627281189, 7.169353723526001, -0.4825034886598586, 0.0609338097274303, 0.0045814891345798, 0.0629191175103187),
    @(-3.050878047943115, 7.165829300880432, -0.3550609424710275, 0.0303905457258224, 0.0102319931611418, 0.0383317954838275),
    @(-3.056754767894745, 7.225212574005127, -0.4760921187698841, -0.0080939643085002, 0.0009162978967650999, -0.0238237436860799),
    @(-3.037489891052246, 7.236634731292725, -0.4997432827949525, 0.0160352122038602, -0.0207694191485643, 0.0113010071218013),
    @(-2.916959762573243, 7.203977525234222, -0.616998553276062, -0.0288633834570646, -0.00137444678694, -0.0123700210824608),
    @(-3.009585857391357, 7.19498348236084, -0.6860059350728989, 0.0059559359215199, 0.0401643887162208, 0.0155770638957619),
    @(-2.86443132162094, 7.121285438537598, -0.5276834592223163, 0.0117591563612222, 0.0131336031481623, 0.0174096599221229),
    @(-2.889585494995118, 7.118069887161255, -0.4351722449064255, 0.0245873257517814, 0.0372627787292003, 0.028557950630784),
    @(-3.056696653366089, 7.102567493915558, -0.6014280728995802, 0.0038179077673703, 0.0119118718430399, -0.0335975885391235)
)

for ($i = 0; $i -lt $newTopData.Count; $i++) {
    $row = 2 + $i
    $vals = $newTopData[$i]
    $ws.Cells.Item($row, 2).Value = "falling"
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
    $ws.Cells.Item($row, 6).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
    $ws.Cells.Item($row, 8).Value = $vals[5]
}

# --- Step 3: append one brand new row (row 31) after the (now shifted) last row ---
$ws.Cells.Item(31, 2).Value = "falling"
$ws.Cells.Item(31, 3).Value = 2.491997003555297
$ws.Cells.Item(31, 4).Value = 6.562706351280213
$ws.Cells.Item(31, 5).Value = -1.225150167942047
$ws.Cells.Item(31, 6).Value = 0.0256563406437635
$ws.Cells.Item(31, 7).Value = 0.0372627787292003
$ws.Cells.Item(31, 8).Value = 0.0050396383740007

# --- Step 4: recompute the contiguous timestamp sequence in column A ---
for ($row = 2; $row -le 31; $row++) {
    $ws.Cells.Item($row, 1).Value = ($row - 2) * 100
}
